{"js": "// Locate the \"Edit Files\" bullet item and add two new bullet items right\n// after it: \"Branch creation \" and \"Merge Branch\".\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet editFilesParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === \"Edit Files\") {\n    editFilesParagraph = paragraphs.items[i];\n    break;\n  }\n}\nif (!editFilesParagraph) {\n  // Fall back to the last paragraph in the body if the text could not be\n  // matched for some reason.\n  editFilesParagraph = paragraphs.items[paragraphs.items.length - 1];\n}\n\nconst branchParagraph = editFilesParagraph.insertParagraph(\n  \"Branch creation \",\n  Word.InsertLocation.after\n);\nawait context.sync();\n\nconst mergeParagraph = branchParagraph.insertParagraph(\n  \"Merge Branch\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"Edit Files\" bullet item (falls back to the document's last\n# paragraph if, for some reason, the text cannot be matched).\n$editFilesParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    $t = $t.TrimEnd(\"`r\", \"`a\", \"`n\")\n    if ($t -eq \"Edit Files\") {\n        $editFilesParagraph = $p\n    }\n}\nif ($editFilesParagraph -eq $null) {\n    $editFilesParagraph = $d.Paragraphs.Item($d.Paragraphs.Count)\n}\n\n# Insert \"Branch creation \" as a new list item right after it.\n$editFilesParagraph.Range.InsertParagraphAfter()\n$count = $d.Paragraphs.Count\n$branchParagraph = $d.Paragraphs.Item($count)\n$branchParagraph.Range.InsertAfter(\"Branch creation \")\n\n# Insert \"Merge Branch\" as a new list item right after that.\n$branchParagraph.Range.InsertParagraphAfter()\n$count = $d.Paragraphs.Count\n$mergeParagraph = $d.Paragraphs.Item($count)\n$mergeParagraph.Range.InsertAfter(\"Merge Branch\")\n"}
